$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD (Wins), AE (Losses), AF (Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, bordered, centered) used by A1:AC1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Team record data for every data row (2-54): Wins=70, Losses=92, Ties=0
$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 70
    $ws.Cells.Item($r, 31).Value = 92
    $ws.Cells.Item($r, 32).Value = 0
}
